{"js": "// 1) \"Time-sensitive\" -> \"Time-saving\"\n//    run text \"ensitive\" is replaced in-place with \"aving\"\n{\n  const results = context.document.body.search(\"ensitive\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for 'ensitive', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"aving\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Post-i\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..services\" -> \"Post-implementation\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..services\"\n//    new text \"mplementation\" inserted immediately after the \"Post-i\" run\n{\n  const results = context.document.body.search(\"Post-i\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for 'Post-i', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"mplementation\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 3) \"Self- sufficient\" -> \"Self-sufficient\" (drop the stray space before the typed \"s\")\n{\n  const results = context.document.body.search(\"Self- s\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for 'Self- s', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"Self-s\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \"Self-e\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026guideline\" -> \"Self-explanatory\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026guideline\"\n//    new text \"xplanatory\" inserted immediately after the \"Self-e\" run\n{\n  const results = context.document.body.search(\"Self-e\", { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for 'Self-e', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"xplanatory\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 5) \"N\u2026\u2026\u2026\u2026\u2026\u2026-core activities\" -> \"Non\u2026\u2026\u2026\u2026\u2026\u2026-core activities\"\n//    new text \"on\" inserted immediately after the leading \"N\"\n{\n  const ellipsis6 = String.fromCharCode(8230).repeat(6);\n  const needle = \"N\" + ellipsis6 + \"-core activities\";\n  const results = context.document.body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for 'N...-core activities', found \" + results.items.length);\n  }\n  const replacement = \"N\" + \"on\" + ellipsis6 + \"-core activities\";\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 6) \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026of communication is inevitable...\" -> \"\u2026Breakdown\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026of communication is inevitable...\"\n//    new text \"Breakdown\" inserted after the first ellipsis character\n{\n  const tail = \"of communication is inevitable when one of the sides is reluctant to compromise\";\n  const ellipsis10 = String.fromCharCode(8230).repeat(10);\n  const needle = ellipsis10 + tail;\n  const results = context.document.body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for the communication sentence, found \" + results.items.length);\n  }\n  const ellipsis1 = String.fromCharCode(8230).repeat(1);\n  const ellipsis9 = String.fromCharCode(8230).repeat(9);\n  const replacement = ellipsis1 + \"Breakdown\" + ellipsis9 + tail;\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) \"Time-sensitive\" -> \"Time-saving\"\n#    The existing \"ensitive\" text is overwritten in place.\n# ---------------------------------------------------------------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"ensitive\")\nif (-not $found) { throw \"Could not find 'ensitive'\" }\n$range.Text = \"aving\"\n\n# ---------------------------------------------------------------------------\n# 2) \"Post-i\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026services\" -> \"Post-implementation\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026services\"\n#    \"mplementation\" is inserted immediately after \"Post-i\".\n# ---------------------------------------------------------------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"Post-i\")\nif (-not $found) { throw \"Could not find 'Post-i'\" }\n$range.Collapse(0)   # wdCollapseEnd\n$range.InsertAfter(\"mplementation\")\n\n# ---------------------------------------------------------------------------\n# 3) \"Self- sufficient\" -> \"Self-sufficient\"\n#    Drop the stray space between \"Self-\" and the typed \"s\".\n# ---------------------------------------------------------------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"Self- s\")\nif (-not $found) { throw \"Could not find 'Self- s'\" }\n$range.Text = \"Self-s\"\n\n# ---------------------------------------------------------------------------\n# 4) \"Self-e\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026guideline\" -> \"Self-explanatory\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026guideline\"\n#    \"xplanatory\" is inserted immediately after \"Self-e\".\n# ---------------------------------------------------------------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"Self-e\")\nif (-not $found) { throw \"Could not find 'Self-e'\" }\n$range.Collapse(0)   # wdCollapseEnd\n$range.InsertAfter(\"xplanatory\")\n\n# ---------------------------------------------------------------------------\n# 5) \"N\u2026\u2026\u2026\u2026\u2026\u2026-core activities\" -> \"Non\u2026\u2026\u2026\u2026\u2026\u2026-core activities\"\n#    \"on\" is inserted immediately after the leading \"N\".\n# ---------------------------------------------------------------------------\n$dots6 = \"\u2026\u2026\u2026\u2026\u2026\u2026\"\n$range = $d.Content\n$needle = \"N\" + $dots6 + \"-core activities\"\n$found = $range.Find.Execute($needle)\nif (-not $found) { throw \"Could not find 'N...-core activities'\" }\n$replacement = \"Non\" + $dots6 + \"-core activities\"\n$range.Text = $replacement\n\n# ---------------------------------------------------------------------------\n# 6) \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026of communication is inevitable...\" -> \"\u2026Breakdown\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026of communication is inevitable...\"\n#    \"Breakdown\" is inserted right after the first ellipsis character.\n# ---------------------------------------------------------------------------\n$tail = \"of communication is inevitable when one of the sides is reluctant to compromise\"\n$dots10 = \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\"\n$dots1 = \"\u2026\"\n$dots9 = \"\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\"\n$range = $d.Content\n$needle = $dots10 + $tail\n$found = $range.Find.Execute($needle)\nif (-not $found) { throw \"Could not find the communication sentence\" }\n$replacement = $dots1 + \"Breakdown\" + $dots9 + $tail\n$range.Text = $replacement\n"}
